$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = 0, with bold font, thin box border, centered horizontally, top vertically
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160

# Reuse the exact same style for A2 via copy/paste-special (formats only),
# so the engine dedupes to a single shared cell style instead of minting a
# near-duplicate one.
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0

# B2 = "disconnected_elements" (plain, unstyled -> shared string)
$ws.Range("B2").Value = "disconnected_elements"
